# framework_offline.pptx - "updates to framework section"
#
# 1. Slide 1 (grouped diagram "Group 15"):
#    - Rounded Rectangle 16 ("Virtual system" -> "Planning system"), and widen it slightly
#    - Rounded Rectangle 17 ("Real system" -> "Tracking system", merged back into one run)
#    - TextBox 26 ("Real dynamics" -> "Tracking dynamics")
#    - TextBox 27 ("Virtual dynamics" -> "Planning dynamics")
# 2. Slide master + every slide layout: the cached "datetimeFigureOut" field text
#    ("17/03/02" -> "3/15/17")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# --- "Virtual system" -> "Planning system" (also gets a touch wider) ---
$virtualSystem = $grp.GroupItems.Item(1)
$virtualSystem.TextFrame.TextRange.Text = "Planning system"
$virtualSystem.Width = 1075467 / 914400 * 72

# --- "Real " + "system" -> single run "Tracking system" ---
$realSystem = $grp.GroupItems.Item(2)
$realSystem.TextFrame.TextRange.Text = ""
$realSystem.TextFrame.TextRange.Text = "Tracking system"

# --- "Real dynamics" -> "Tracking dynamics" ---
$realDynamics = $grp.GroupItems.Item(8)
$realDynamics.TextFrame.TextRange.Text = "Tracking dynamics"

# --- "Virtual dynamics" -> "Planning dynamics" ---
$virtualDynamics = $grp.GroupItems.Item(9)
$virtualDynamics.TextFrame.TextRange.Text = "Planning dynamics"

# --- Refresh the cached date field text everywhere it is cached ---
function Update-CachedDate($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "17/03/02") {
                $sh.TextFrame.TextRange.Text = "3/15/17"
            }
        }
    }
}

Update-CachedDate $p.SlideMaster.Shapes
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Update-CachedDate $layout.Shapes
}
